$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5185.875
$ws.Range("I40").Value = 2894.5
$ws.Range("J40").Value = 5949.6665
$ws.Range("K40").Value = 2894.5
$ws.Range("L40").Value = 5949.6665
$ws.Range("M40").Value = -2719.5
$ws.Range("N40").Value = -6299.6665

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H44").Value = 37745
$ws.Range("J44").Value = 37745
$ws.Range("L44").Value = 37745
$ws.Range("N44").Value = -38669

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 11513.637
$ws.Range("I62").Value = 14258.25
$ws.Range("K62").Value = 14258.25
$ws.Range("M62").Value = -13634.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 11513.637
$ws.Range("I65").Value = 14258.25
$ws.Range("K65").Value = 71291.25
$ws.Range("M65").Value = -68171.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 7436.2856
$ws.Range("I132").Value = 6766.615
$ws.Range("K132").Value = 20299.845
$ws.Range("M132").Value = -17769.845

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 59263.65
$ws.Range("I135").Value = 465.375
$ws.Range("J135").Value = 1000036
$ws.Range("K135").Value = 4188.375
$ws.Range("L135").Value = 9000324
$ws.Range("M135").Value = -1653.375
$ws.Range("N135").Value = -9005394

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2992.578
$ws.Range("J138").Value = 3023.5
$ws.Range("L138").Value = 9070.5
$ws.Range("N138").Value = -19350.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 29257.79
$ws.Range("I32").Value = 31312.031
$ws.Range("K32").Value = 31312.031
$ws.Range("M32").Value = -31025.031

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H50").Value = 13484.75
$ws.Range("J50").Value = 14979.667
$ws.Range("L50").Value = 14979.667
$ws.Range("N50").Value = -16407.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 37999
$ws.Range("J76").Value = 37999
$ws.Range("L76").Value = 37999
$ws.Range("N76").Value = -38675

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H79").Value = 37999
$ws.Range("J79").Value = 37999
$ws.Range("L79").Value = 37999
$ws.Range("N79").Value = -40339

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H86").Value = 22795
$ws.Range("I86").Value = 22795
$ws.Range("K86").Value = 22795
$ws.Range("M86").Value = -21609

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H89").Value = 22795
$ws.Range("I89").Value = 22795
$ws.Range("K89").Value = 68385
$ws.Range("M89").Value = -62457

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H121").Value = 59133.332
$ws.Range("J121").Value = 59133.332
$ws.Range("L121").Value = 59133.332
$ws.Range("N121").Value = -62627.332

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 5443.4287
$ws.Range("I122").Value = 4984
$ws.Range("K122").Value = 14952
$ws.Range("M122").Value = -12502

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 9095381
$ws.Range("I20").Value = 20001160
$ws.Range("J20").Value = 7232.1665
$ws.Range("K20").Value = 20001160
$ws.Range("L20").Value = 7232.1665
$ws.Range("M20").Value = -20000913
$ws.Range("N20").Value = -7726.1665

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H119").Value = 250000
$ws.Range("J119").Value = 250000
$ws.Range("L119").Value = 250000
$ws.Range("N119").Value = -259676

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1496.125
$ws.Range("I16").Value = 1468.6923
$ws.Range("J16").Value = 1615
$ws.Range("K16").Value = 1468.6923
$ws.Range("L16").Value = 1615
$ws.Range("M16").Value = -1181.6923
$ws.Range("N16").Value = -2189

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4302.6
$ws.Range("I31").Value = 2218.077
$ws.Range("K31").Value = 2218.077
$ws.Range("M31").Value = -1923.077

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4302.6
$ws.Range("I34").Value = 2218.077
$ws.Range("K34").Value = 2218.077
$ws.Range("M34").Value = -2016.077

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1496.125
$ws.Range("I113").Value = 1468.6923
$ws.Range("J113").Value = 1615
$ws.Range("K113").Value = 1468.6923
$ws.Range("L113").Value = 1615
$ws.Range("M113").Value = 701.3077000000001
$ws.Range("N113").Value = -5955

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 375.7857
$ws.Range("I107").Value = 233
$ws.Range("J107").Value = 482.875
$ws.Range("K107").Value = 699
$ws.Range("L107").Value = 1448.625
$ws.Range("M107").Value = 1221
$ws.Range("N107").Value = -5288.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3424.1538
$ws.Range("I70").Value = 2131
$ws.Range("J70").Value = 4932.8335
$ws.Range("K70").Value = 2131
$ws.Range("L70").Value = 4932.8335
$ws.Range("M70").Value = -1861
$ws.Range("N70").Value = -5472.8335

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 3424.1538
$ws.Range("I73").Value = 2131
$ws.Range("J73").Value = 4932.8335
$ws.Range("K73").Value = 2131
$ws.Range("L73").Value = 4932.8335
$ws.Range("M73").Value = -1195
$ws.Range("N73").Value = -6804.8335

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3417.1304
$ws.Range("I102").Value = 2533.111
$ws.Range("J102").Value = 6599.6
$ws.Range("K102").Value = 2533.111
$ws.Range("L102").Value = 6599.6
$ws.Range("M102").Value = -911.1109999999999
$ws.Range("N102").Value = -9843.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1817.9474
$ws.Range("I22").Value = 1022.7143
$ws.Range("K22").Value = 1022.7143
$ws.Range("M22").Value = -727.7143

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1817.9474
$ws.Range("I27").Value = 1022.7143
$ws.Range("K27").Value = 1022.7143
$ws.Range("M27").Value = -915.7143

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 6623.6924
$ws.Range("I61").Value = 5237
$ws.Range("K61").Value = 5237
$ws.Range("M61").Value = -5035

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1656
$ws.Range("I93").Value = 1651.1177
$ws.Range("K93").Value = 1651.1177
$ws.Range("M93").Value = -403.1177

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 6623.6924
$ws.Range("I113").Value = 5237
$ws.Range("K113").Value = 5237
$ws.Range("M113").Value = -3067

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 78457.44
$ws.Range("I132").Value = 102235.586
$ws.Range("K132").Value = 306706.758
$ws.Range("M132").Value = -304176.758

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 6460.5386
$ws.Range("I136").Value = 5787.5557
$ws.Range("K136").Value = 17362.6671
$ws.Range("M136").Value = -14812.6671

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2175.7097
$ws.Range("I107").Value = 1005.75
$ws.Range("K107").Value = 3017.25
$ws.Range("M107").Value = -1097.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 64340.188
$ws.Range("I126").Value = 64340.188
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 193020.564
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -190550.564
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H127").Value = 98452
$ws.Range("J127").Value = 98452
$ws.Range("L127").Value = 98452
$ws.Range("N127").Value = -108372

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 74402
$ws.Range("I132").Value = 93977.37
$ws.Range("J132").Value = 2625.6667
$ws.Range("K132").Value = 281932.11
$ws.Range("L132").Value = 7877.000100000001
$ws.Range("M132").Value = -279402.11
$ws.Range("N132").Value = -12937.0001
